$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B3").Value = 6.143100000000001
$ws.Range("D3").Value = -6.990299999999994
$ws.Range("B21").Value = 9.373600000000005
$ws.Range("B23").Value = 9.180200000000008
$ws.Range("D24").Value = -7.525800000000004
$ws.Range("B25").Value = 6.318399999999998
$ws.Range("C27").Value = -13.1133
$ws.Range("C31").Value = -13.3039
$ws.Range("C39").Value = -12.52050000000001
$ws.Range("C48").Value = -11.4033
$ws.Range("C51").Value = -11.5478
$ws.Range("C52").Value = -11.11369999999999
$ws.Range("B53").Value = 5.149800000000002
$ws.Range("C55").Value = -13.7441
$ws.Range("C56").Value = -12.86819999999999
$ws.Range("B57").Value = 4.853399999999995
$ws.Range("C57").Value = -13.50789999999999
$ws.Range("D57").Value = -8.639799999999994
$ws.Range("B59").Value = 4.712899999999995
$ws.Range("D61").Value = -7.807899999999998
$ws.Range("B69").Value = 5.432999999999994
$ws.Range("D70").Value = -7.244899999999997
$ws.Range("C73").Value = -12.22600000000001
$ws.Range("B79").Value = 9.104200000000002
$ws.Range("B83").Value = 5.605
$ws.Range("D86").Value = -7.698899999999994
$ws.Range("C89").Value = -11.089
$ws.Range("C90").Value = -12.337
$ws.Range("B93").Value = 5.519499999999998
$ws.Range("D98").Value = -8.608099999999997
$ws.Range("D100").Value = -8.156700000000003
$ws.Range("D102").Value = -7.862799999999996
